$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The old single "wait_reward_coef" row (row 7: wait_reward_coef /
# "Balancing coefficient for wait in reward calculation" / 1) is split into
# two separate parameters - one for vehicles, one for pedestrians:
#   row 7 -> wait_veh_reward_coef  (alpha_veh ...)
#   row 8 -> wait_ped_reward_coef  (alpha_ped ...)
# Everything that used to live at rows 8..27 shifts down to rows 9..28.
# ---------------------------------------------------------------------------

# Insert a new row at 8. This shifts old rows 8..27 down to 9..28 and keeps
# the old row 7 ("wait_reward_coef" / "Balancing coefficient...") at row 7
# for the moment - we overwrite it below.
$ws.Rows(8).Insert()

# The freshly inserted row 8 doesn't reliably inherit the "D" column number
# style (border + centred Consolas) used by every sibling parameter row, so
# copy it explicitly from a neighbouring row before filling in values.
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 7: wait_veh_reward_coef -------------------------------------------
$ws.Range("B7").Value2 = "wait_veh_reward_coef"

$ws.Range("C7").Value2 = "α_veh Balancing coefficient for wait_veh in reward calculation"
$ws.Range("C7").WrapText = $true
$ws.Range("C7").VerticalAlignment = -4108
$ws.Range("C7").Font.Name = "Consolas"
$ws.Range("C7").Font.Size = 11
$ws.Range("C7").Font.Color = 0
$cVeh = $ws.Range("C7").Characters(6, 57)
$cVeh.Font.Name = "Consolas"
$cVeh.Font.Size = 11
$cVeh.Font.Color = 0

$ws.Range("D7").Value2 = 1

$ws.Rows(7).RowHeight = 28.8

# --- Row 8: wait_ped_reward_coef -------------------------------------------
$ws.Range("B8").Value2 = "wait_ped_reward_coef"

$ws.Range("C8").Value2 = "α_ped Balancing coefficient for wait_ped in reward calculation"
$ws.Range("C8").WrapText = $true
$ws.Range("C8").VerticalAlignment = -4108
$ws.Range("C8").Font.Name = "Consolas"
$ws.Range("C8").Font.Size = 11
$ws.Range("C8").Font.Color = 0
$cPed = $ws.Range("C8").Characters(6, 57)
$cPed.Font.Name = "Consolas"
$cPed.Font.Size = 11
$cPed.Font.Color = 0

$ws.Range("D8").Value2 = 1

$ws.Rows(8).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Restore the selected cell in the frozen bottom-right pane.
# ---------------------------------------------------------------------------
$ws.Range("C11").Select()
